$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers but must remain text
# (matching the original inlineStr/text cell type in the workbook).
# We force the NumberFormat to Text, assign the value, then clear the
# format again so the cell keeps its original (default) style.
$textCells = @(
    @{Addr='D5'; Val='320.28'}
    @{Addr='D7'; Val='0.5056'}
    @{Addr='D8'; Val='0.4072'}
    @{Addr='D9'; Val='0.08347'}
    @{Addr='D10'; Val='1.110'}
    @{Addr='D11'; Val='42.31'}
    @{Addr='D12'; Val='24.06'}
    @{Addr='D13'; Val='6.424'}
    @{Addr='D15'; Val='7.243'}
    @{Addr='D17'; Val='92.54'}
    @{Addr='D19'; Val='0.06505'}
    @{Addr='D20'; Val='18.51'}
    @{Addr='D22'; Val='5.942'}
    @{Addr='D25'; Val='2.194'}
    @{Addr='D28'; Val='162.67'}
    @{Addr='D30'; Val='128.92'}
    @{Addr='D31'; Val='1.135'}
    @{Addr='D34'; Val='3.786'}
    @{Addr='D35'; Val='0.02456'}
    @{Addr='D36'; Val='5.339'}
    @{Addr='D38'; Val='0.2155'}
    @{Addr='D39'; Val='0.6516'}
    @{Addr='D40'; Val='1.197'}
    @{Addr='D41'; Val='8.620'}
    @{Addr='D43'; Val='1.209'}
    @{Addr='D44'; Val='13.47'}
    @{Addr='D45'; Val='0.6067'}
    @{Addr='D46'; Val='2.187'}
    @{Addr='D48'; Val='1.211'}
    @{Addr='D49'; Val='122.21'}
    @{Addr='D50'; Val='79.00'}
    @{Addr='D51'; Val='1.142'}
)

foreach ($item in $textCells) {
    $rng = $ws.Range($item.Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Val
    $rng.ClearFormats()
}

# Remaining cells: plain text assignment is sufficient since Excel
# will not reinterpret these strings as numbers.
$ws.Range("D2").Value = '30.118.20'
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").Value = '1.917.67'
$ws.Range("E3").Value = '  +2.53%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("E5").Value = '  +0.24%  '
$ws.Range("E7").Value = '  -0.81%  '
$ws.Range("E8").Value = '  +3.53%  '
$ws.Range("E9").Value = '  +1.75%  '
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("E10").Value = '  +1.57%  '
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("E11").Value = '  +0.18%  '
$ws.Range("E12").Value = '  +4.96%  '
$ws.Range("E13").Value = '  +2.43%  '
$ws.Range("D14").Value = '1.914.39'
$ws.Range("E14").Value = '  +2.90%  '
$ws.Range("E15").Value = '  +1.20%  '
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("E17").Value = '  +0.74%  '
$ws.Range("E18").Value = '  +1.28%  '
$ws.Range("E19").Value = '  +1.37%  '
$ws.Range("E20").Value = '  +3.03%  '
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("E22").Value = '  +2.30%  '
$ws.Range("D23").Value = '30.127.52'
$ws.Range("E23").Value = '  +0.39%  '
$ws.Range("E24").Value = '  +2.46%  '
$ws.Range("E25").Value = '  +1.36%  '
$ws.Range("D26").Value = '2.134.06'
$ws.Range("E26").Value = '  +2.69%  '
$ws.Range("E27").Value = '  +4.34%  '
$ws.Range("E28").Value = '  +0.94%  '
$ws.Range("E29").Value = '  +1.70%  '
$ws.Range("E30").Value = '  +1.45%  '
$ws.Range("E31").Value = '  +7.59%  '
$ws.Range("E32").Value = '  +0.95%  '
$ws.Range("E33").Value = '  +0.75%  '
$ws.Range("E34").Value = '  +1.09%  '
$ws.Range("E35").Value = '  +1.63%  '
$ws.Range("E36").Value = '  +1.15%  '
$ws.Range("E37").Value = '  +1.81%  '
$ws.Range("E38").Value = '  +0.54%  '
$ws.Range("E39").Value = '  +3.38%  '
$ws.Range("E40").Value = '  +2.05%  '
$ws.Range("E41").Value = '  +1.27%  '
$ws.Range("E42").Value = '  +1.36%  '
$ws.Range("E43").Value = '  +0.57%  '
$ws.Range("E44").Value = '  +4.70%  '
$ws.Range("E45").Value = '  +2.87%  '
$ws.Range("E46").Value = '  +9.79%  '
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("B48").Value = 'EOS'
$ws.Range("C48").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("E48").Value = '  +0.70%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("E49").Value = '  -0.30%  '
$ws.Range("E50").Value = '  +3.19%  '
$ws.Range("E51").Value = '  +0.10%  '
